# Auto-generated: apply cryptos price/volume/coin updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.408.11"
$ws.Range("E2").Value = "  +4.82%  "
$ws.Range("D3").Value = "3.604.38"
$ws.Range("E3").Value = "  +9.48%  "
$ws.Range("D5").Value = "'240.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.82%  "
$ws.Range("D6").Value = "'638.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("E7").Value = "  +9.51%  "
$ws.Range("D8").Value = "'0.404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.11%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +10.12%  "
$ws.Range("D11").Value = "3.600.96"
$ws.Range("E11").Value = "  +9.50%  "
$ws.Range("D12").Value = "'43.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.58%  "
$ws.Range("E13").Value = "  +4.94%  "
$ws.Range("D14").Value = "'6.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.60%  "
$ws.Range("D15").Value = "4.291.33"
$ws.Range("E15").Value = "  +9.97%  "
$ws.Range("D16").Value = "96.299.52"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("E17").Value = "  +6.10%  "
$ws.Range("D18").Value = "'8.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.16%  "
$ws.Range("D19").Value = "3.596.66"
$ws.Range("E19").Value = "  +9.28%  "
$ws.Range("D20").Value = "'13.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +25.27%  "
$ws.Range("D21").Value = "'18.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.87%  "
$ws.Range("D22").Value = "'0.501"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.47%  "
$ws.Range("D23").Value = "'516.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.99%  "
$ws.Range("D24").Value = "'3.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("D25").Value = "'0.0000199"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.68%  "
$ws.Range("D26").Value = "'6.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.49%  "
$ws.Range("D27").Value = "'97.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.58%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'12.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.61%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'3.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +19.13%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'11.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.53%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.143"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.01%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").Value = "'0.182"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.15%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'30.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.21%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "'0.569"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.22%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'575.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.88%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'7.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.38%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'1.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.69%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.152"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.45%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.929"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.55%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").Value = "'1.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.60%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0432"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.25%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "'23.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'5.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.59%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.24%  "
$ws.Range("D48").Value = "'3.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'53.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.86%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'8.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'3.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.04%  "
